$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the most recent weekly sheet (Nov26Data) to create the new
#    week's sheet (Dec03Data), placed at the end of the tab strip.
# ---------------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("Nov26Data")
$lastSheet   = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Dec03Data"

# ---------------------------------------------------------------------------
# 2. Write this week's fantasy point data into the new sheet.
#    columns: name | position_1 | position_2 | team | games_7 | points_7 |
#             games_total | points_total | games_this_week
# ---------------------------------------------------------------------------
$rows = @(
    @("Marner",      "C", "R",  "TOR", 3, 44.7,               27, 288.10000000000002, 3),
    @("Skinner",     "C", "L",  "BUF", 3, 18.7,               27, 280.7,              3),
    @("Wheeler",     "C", "R",  "WPG", 4, 47.8,               26, 267.8,              3),
    @("Backstrom",   "C", "",   "WAS", 3, 63.1,               26, 264,                3),
    @("Kane",        "R", "",   "CHI", 4, 10.199999999999999, 27, 244.6,              3),
    @("Monahan",     "C", "",   "CGY", 3, 25.7,               27, 242.7,              4),
    @("Hall",        "L", "",   "NJD", 3, 26.9,               25, 234.1,              3),
    @("Parise",      "L", "",   "MIN", 3, 26.1,               25, 226.8,              3),
    @("Pettersson",  "C", "",   "VAN", 3, 6.7,                23, 196.1,              3),
    @("Suter",       "D", "",   "MIN", 3, 20.7,               26, 178,                3),
    @("Ekholm",      "D", "",   "NSH", 3, 18.600000000000001, 27, 177.9,              3),
    @("Pionk",       "D", "",   "NYR", 4, 18.399999999999999, 26, 176.6,              1),
    @("Ellis",       "D", "",   "NSH", 3, 10.3,               27, 171.5,              3),
    @("Seabrook",    "D", "",   "CHI", 4, 8.6999999999999993, 28, 168.4,              3),
    @("Hellebuyck",  "G", "",   "WPG", 3, 14.8,               20, 202.6,              2),
    @("Bishop",      "G", "",   "DAL", 1, 16.399999999999999, 16, 207.4,              2)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq "") {
        $newSheet.Cells.Item($r, 3).ClearContents()
    } else {
        $newSheet.Cells.Item($r, 3).Value = $row[2]
    }
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $newSheet.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# points_7 / points_total columns get a 2-decimal number format.
$newSheet.Range("F2:F17").NumberFormat = "0.00"
$newSheet.Range("H2:H17").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 3. Update leftover view/selection state on the sheets that were touched
#    while building the new sheet.
# ---------------------------------------------------------------------------
# Nov26Data had its header row selected (e.g. while copying it to build the
# new sheet) before the new sheet took over as the active tab.
$sourceSheet.Activate()
$sourceSheet.Rows("1:1").Select()

# The new sheet ends up selected with its 9th row highlighted, and becomes
# the active tab shown when the workbook is reopened.
$newSheet.Activate()
$newSheet.Rows("9:9").Select()
